$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting existing rows 154..274 down to 155..275.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new weekly price record.
$ws.Cells.Item(154, 1).Value  = 3
$ws.Cells.Item(154, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(154, 3).Value  = "Coquimbo"
$ws.Cells.Item(154, 4).Value  = 44512
$ws.Cells.Item(154, 5).Value  = 5
$ws.Cells.Item(154, 6).Value  = 100112028
$ws.Cells.Item(154, 7).Value  = "Sandia"
$ws.Cells.Item(154, 8).Value  = "Sin especificar"
$ws.Cells.Item(154, 9).Value  = "Primera"
$ws.Cells.Item(154, 10).Value = 180
$ws.Cells.Item(154, 11).Value = 600
$ws.Cells.Item(154, 12).Value = 600
$ws.Cells.Item(154, 13).Value = 600
$ws.Cells.Item(154, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(154, 15).Value = "Perú"
$ws.Cells.Item(154, 16).Value = 600
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
